# Reorder the comma-separated "Recorded By" names in column G so that
# multi-valued entries are listed in a different (swapped) order, matching
# the latest sync from the source system.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact lookup table of how each distinct multi-value "Recorded By" string
# should be reordered. Single-value cells are left untouched.
$map = @{
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"           = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"          = "System, backup@backdoor.com"
    "admin@admin.com, System"              = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com"  = "dnasr281@gmail.com, admin@admin.com"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
